$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): "web link" -> "web link/data source" ----------
$ws.Range("D1").Value = "web link/data source"

# --- Two new rows documenting the tables downloaded from proteinatlas.org
$ws.Cells.Item(7, 1).Value = "HLA_aggregated.tsv"
$ws.Cells.Item(7, 4).Value = "https://www.proteinatlas.org/"
$ws.Cells.Item(8, 1).Value = "proteinatlas.tsv"
$ws.Cells.Item(8, 4).Value = "https://www.proteinatlas.org/"

# --- Restyle the header row: drop from 16pt to 12pt (keep bold), and let
#     the row height return to the sheet default instead of a fixed 21.
$ws.Range("A1:D1").Font.Size = 12
$ws.Rows.Item(1).AutoFit()

# --- Restyle row 2 ("dash.csv"): drop from 14pt to the default 12pt, and
#     shrink the (now smaller) wrapped row from 60 to 51.
$ws.Range("A2:D2").Font.Size = 12
$ws.Rows.Item(2).RowHeight = 51

# --- New "web link/data source" cells get the same muted grey used by
#     modern link styling (#212529) at 12pt.
$ws.Range("D7:D8").Font.Size = 12
$ws.Range("D7:D8").Font.Color = 2696481

# --- Restore the whole-sheet selection state recorded for this edit.
$ws.Cells.Select() | Out-Null
